# Apply the edit described by the diff to LuanMenh.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New rows (42-54) appended to Sheet2, columns A and B.
# Each pair is (A value, B value)
$rows = @(
    @("Tử Vi tọa thủ cung Mệnh tại Thìn", "Tử Vi tọa thủ cung Mệnh tại Thìn"),
    @("Phá Quân tọa thủ cung Mệnh tại Thìn", "Phá Quân tọa thủ cung Mệnh tại Thìn"),
    @("Tử Vi tọa thủ cung Mệnh tại Tuất", "Phá Quân tọa thủ cung Mệnh tại Tuất"),
    @("Tử Vi tọa thủ cung Mệnh tại Mùi", "Phá Quân tọa thủ cung Mệnh tại Mùi"),
    @("Tử Vi tọa thủ cung Mệnh tại Sửu", "Phá Quân tọa thủ cung Mệnh tại Sửu"),
    @("Tử Vi tọa thủ cung Mệnh tại Thìn gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Tử Vi tọa thủ cung Mệnh tại Thìn gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Phá Quân tọa thủ cung Mệnh tại Thìn gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Phá Quân tọa thủ cung Mệnh tại Thìn gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Tử Vi tọa thủ cung Mệnh tại Tuất gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Tử Vi tọa thủ cung Mệnh tại Tuất gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Phá Quân tọa thủ cung Mệnh tại Tuất gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Phá Quân tọa thủ cung Mệnh tại Tuất gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Tử Vi tọa thủ cung Mệnh tại Mùi gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Tử Vi tọa thủ cung Mệnh tại Mùi gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Phá Quân tọa thủ cung Mệnh tại Mùi gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Phá Quân tọa thủ cung Mệnh tại Mùi gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Tử Vi tọa thủ cung Mệnh tại Sửu gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Tử Vi tọa thủ cung Mệnh tại Sửu gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc"),
    @("Phá Quân tọa thủ cung Mệnh tại Sửu gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc", "Phá Quân tọa thủ cung Mệnh tại Sửu gặp các sao cát tinh: Hóa Quyền, Hóa Lộc, Hóa Khoa, Thiên Phủ, Tả Phù, Hữu Bật, Thiên Tướng, Văn Xương, Văn Khúc")
)

$startRow = 42
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Update the view so it mirrors the author's final view state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O45").Select()
